# Update workspace configuration and add LeetCode problems 53 (Maximum
# Subarray), 2322 (Minimum Score After Removals on a Tree), 621 (Task
# Scheduler) and 233 (Number of Digit One) to the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook window geometry (best effort) ---------------------------
try {
    $win = $wb.Windows.Item(1)
    $win.Width = 29100
} catch {
}

# --- Template cells used to clone the date-column number format -------
$dateTemplate = $ws.Range("I92")

# --- Row 93: Maximum Subarray ------------------------------------------
$ws.Rows(93).RowHeight = 68
$ws.Range("A93").Value2 = 53
$ws.Range("B93").Value2 = "Maximum Subarray"
$ws.Range("C93").Value2 = "#array #divide-and-conquer #dynamic-programming "
$ws.Range("D93").Value2 = "medium"
$ws.Range("E93").Value2 = 1
$ws.Range("F93").Value2 = 1
$ws.Range("G93").Value2 = 20
$ws.Range("I93").Value2 = 45861
$dateTemplate.Copy()
$ws.Range("I93").PasteSpecial(-4122)

# --- Row 94: Minimum Score After Removals on a Tree --------------------
$ws.Rows(94).RowHeight = 51
$ws.Range("A94").Value2 = 2322
$ws.Range("B94").Value2 = "Minimum Score After Removals on a Tree"
$ws.Range("C94").Value2 = "#array #bit-manipulation #tree #dfs"
$ws.Range("D94").Value2 = "hard"
$ws.Range("E94").Value2 = 0
$ws.Range("F94").Value2 = 1
$ws.Range("G94").Value2 = 65
$ws.Range("H94").Value2 = 45862
$ws.Range("I94").Value2 = 45862
$ws.Range("J94").Value2 = "?"
$dateTemplate.Copy()
$ws.Range("H94:I94").PasteSpecial(-4122)

# --- Row 95: Task Scheduler --------------------------------------------
$ws.Rows(95).RowHeight = 34
$ws.Range("A95").Value2 = 621
$ws.Range("B95").Value2 = "Task Scheduler"
$ws.Range("C95").Value2 = "#array #greedy #queue"
$ws.Range("D95").Value2 = "medium"
$ws.Range("E95").Value2 = 1
$ws.Range("F95").Value2 = 1
$ws.Range("G95").Value2 = 50
$ws.Range("H95").Value2 = 45862
$ws.Range("I95").Value2 = 45862
$ws.Range("J95").Value2 = "?"
$dateTemplate.Copy()
$ws.Range("H95:I95").PasteSpecial(-4122)

# --- Row 96: Number of Digit One ----------------------------------------
$ws.Rows(96).RowHeight = 17
$ws.Range("A96").Value2 = 233
$ws.Range("B96").Value2 = "Number of Digit One"
$ws.Range("C96").Value2 = "#math"
$ws.Range("D96").Value2 = "hard"
$ws.Range("E96").Value2 = 0
$ws.Range("F96").Value2 = 1
$ws.Range("G96").Value2 = 40
$ws.Range("H96").Value2 = 45863
$ws.Range("I96").Value2 = 45863
$ws.Range("J96").Value2 = "?"
$dateTemplate.Copy()
$ws.Range("H96:I96").PasteSpecial(-4122)

# --- Selection / viewport (best effort) ---------------------------------
$ws.Range("J95").Select()
try {
    $excel.ActiveWindow.ScrollRow = 89
    $excel.ActiveWindow.ScrollColumn = 2
} catch {
}
